$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The worksheet is protected; temporarily unprotect so the cells below can
# be written, then re-protect it afterwards so the sheet's protected state
# is preserved.
$ws.Unprotect()

# Update the "as of" date in the confidential banner text (A9).
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-11 for illustrative purposes only and are subject to change."
# Setting a wrapped multi-line value auto-marks the row with an explicit
# custom height; AutoFit it back so the row keeps using the default height
# (matching the original, which has no explicit row height on row 9).
$ws.Rows.Item(9).AutoFit()

# Refresh the Weight / Percent Change figures for each sector row.
$ws.Range("D2").Value = 0.255496338694521
$ws.Range("E2").Value = 0.003315043438500087

$ws.Range("D3").Value = 0.2544529250316451
$ws.Range("E3").Value = -0.01774834437086104

$ws.Range("D4").Value = 0.2469107814524077
$ws.Range("E4").Value = -0.01479374110953058

$ws.Range("D5").Value = 0.2431399548214262
$ws.Range("E5").Value = -0.002865329512893977

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = -0.008018546944756788

$ws.Protect()
